# Fix date typo in CV: the "Festive Colleague, Tesco" role shows an end
# date of 12/2024, but it should read 12/2023 (matching the 11/2023 start
# date / a one-month festive-temp role). Target the specific substring
# "12/2024)" so only the erroneous year inside that parenthetical gets
# corrected, leaving everything else in the document untouched.

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "12/2024)",   # FindText
    $true,        # MatchCase
    $false,       # MatchWholeWord
    $false,       # MatchWildcards
    $false,       # MatchSoundsLike
    $false,       # MatchAllWordForms
    $true,        # Forward
    1,            # Wrap (wdFindContinue)
    $false,       # Format
    "12/2023)",   # ReplaceWith
    2             # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Could not find the '12/2024)' text to fix the date typo."
}

Write-Host "Date typo fixed: $found"
